$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 154, shifting existing rows 154-186 down to 155-187
$ws.Rows.Item(154).Insert()

# Populate the new row 154 with the new record's data
$ws.Cells.Item(154, 1).Value = 10
$ws.Cells.Item(154, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(154, 3).Value = "La Araucanía"
$ws.Cells.Item(154, 4).Value = 44798
$ws.Cells.Item(154, 5).Value = 9
$ws.Cells.Item(154, 6).Value = 100114007
$ws.Cells.Item(154, 7).Value = "Jengibre"
$ws.Cells.Item(154, 8).Value = "Sin especificar"
$ws.Cells.Item(154, 9).Value = "Primera"
$ws.Cells.Item(154, 10).Value = 190
$ws.Cells.Item(154, 11).Value = 17000
$ws.Cells.Item(154, 12).Value = 20000
$ws.Cells.Item(154, 13).Value = 18579
$ws.Cells.Item(154, 14).Value = '$/caja 13 kilos'
$ws.Cells.Item(154, 15).Value = "Perú"
$ws.Cells.Item(154, 16).Value = 1429
$ws.Cells.Item(154, 17).Value = 13
$ws.Cells.Item(154, 18).Value = "Hortaliza"
